# Auto-generated edit script applying the cryptos.xlsx data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.322.79"
$ws.Range("E2").Value = "  -1.98%  "
$ws.Range("D3").Value = "1.792.00"
$ws.Range("E3").Value = "  -1.81%  "
$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "'1.004"
$ws.Range("E5").Value = "  -0.11%  "
$ws.Range("D6").Value = "'306.57"
$ws.Range("E6").Value = "  -1.33%  "
$ws.Range("D7").Value = "'0.4508"
$ws.Range("E7").Value = "  -1.32%  "
$ws.Range("D8").Value = "'0.3594"
$ws.Range("E8").Value = "  -2.28%  "
$ws.Range("D9").Value = "'46.34"
$ws.Range("E9").Value = "  +1.34%  "
$ws.Range("D10").Value = "'0.07074"
$ws.Range("E10").Value = "  -1.09%  "
$ws.Range("D11").Value = "'0.8839"
$ws.Range("E11").Value = "  +1.29%  "
$ws.Range("D12").Value = "'0.07734"
$ws.Range("E12").Value = "  -0.38%  "
$ws.Range("D13").Value = "'19.45"
$ws.Range("E13").Value = "  -0.51%  "
$ws.Range("D14").Value = "1.776.13"
$ws.Range("E14").Value = "  -2.51%  "
$ws.Range("D15").Value = "'5.278"
$ws.Range("E15").Value = "  -0.61%  "
$ws.Range("D16").Value = "'6.317"
$ws.Range("E16").Value = "  -0.87%  "
$ws.Range("D17").Value = "'84.88"
$ws.Range("E17").Value = "  -2.27%  "
$ws.Range("D18").Value = "'1.006"
$ws.Range("E18").Value = "  +0.00%  "
$ws.Range("D19").Value = "'0.000008504"
$ws.Range("E19").Value = "  -2.26%  "
$ws.Range("D20").Value = "'1.004"
$ws.Range("E20").Value = "  -0.12%  "
$ws.Range("D21").Value = "'14.24"
$ws.Range("E21").Value = "  -1.43%  "
$ws.Range("D22").Value = "26.360.05"
$ws.Range("E22").Value = "  -1.99%  "
$ws.Range("D23").Value = "'4.966"
$ws.Range("E23").Value = "  -0.57%  "
$ws.Range("D24").Value = "'10.54"
$ws.Range("E24").Value = "  +0.97%  "
$ws.Range("D25").Value = "1.985.53"
$ws.Range("E25").Value = "  -3.44%  "
$ws.Range("D26").Value = "'1.966"
$ws.Range("E26").Value = "  -1.65%  "
$ws.Range("D27").Value = "'150.97"
$ws.Range("E27").Value = "  -0.26%  "
$ws.Range("D28").Value = "'17.80"
$ws.Range("E28").Value = "  -1.73%  "
$ws.Range("D29").Value = "'2.029"
$ws.Range("E29").Value = "  +4.07%  "
$ws.Range("D30").Value = "'111.64"
$ws.Range("E30").Value = "  -1.64%  "
$ws.Range("D31").Value = "'4.837"
$ws.Range("E31").Value = "  -1.28%  "
$ws.Range("D32").Value = "'0.08678"
$ws.Range("E32").Value = "  -1.22%  "
$ws.Range("D33").Value = "'3.070"
$ws.Range("E33").Value = "  +1.12%  "
$ws.Range("D34").Value = "'2.763"
$ws.Range("E34").Value = "  +8.98%  "
$ws.Range("D35").Value = "'4.439"
$ws.Range("E35").Value = "  -0.80%  "
$ws.Range("D36").Value = "'0.7221"
$ws.Range("E36").Value = "  -3.53%  "
$ws.Range("D37").Value = "'1.101"
$ws.Range("E37").Value = "  -2.65%  "
$ws.Range("D38").Value = "'1.003"
$ws.Range("E38").Value = "  +0.08%  "
$ws.Range("D40").Value = "'0.01925"
$ws.Range("E40").Value = "  -0.73%  "
$ws.Range("D41").Value = "'0.05093"
$ws.Range("E41").Value = "  -0.51%  "
$ws.Range("D42").Value = "'2.857"
$ws.Range("E42").Value = "  -1.74%  "
$ws.Range("D43").Value = "'0.5050"
$ws.Range("E43").Value = "  +1.78%  "
$ws.Range("D44").Value = "'6.814"
$ws.Range("E44").Value = "  -1.41%  "
$ws.Range("D45").Value = "'0.1517"
$ws.Range("E45").Value = "  -4.66%  "
$ws.Range("D46").Value = "'7.997"
$ws.Range("E46").Value = "  -3.50%  "
$ws.Range("D47").Value = "'1.004"
$ws.Range("E47").Value = "  -0.12%  "
$ws.Range("D48").Value = "'0.4626"
$ws.Range("E48").Value = "  -1.16%  "
$ws.Range("D51").Value = "'1.567"
$ws.Range("E51").Value = "  -2.52%  "
$ws.Range("E39").Value = "  -1.24%  "

# Rows 49 and 50 swap coin identity (EnergySwap <-> Quant) with refreshed values
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").Value = "'100.56"
$ws.Range("E49").Value = "  -0.89%  "

$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'9.829"
$ws.Range("E50").Value = "  -2.28%  "
